$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 42611.884594907409
$ws.Range("B7").Value = -10
$ws.Range("C7").Value = 46
$ws.Range("D7").Value = 51
$ws.Range("E7").Value = 40
$ws.Range("F7").Value = 60
$ws.Range("G7").Value = 8104
$ws.Range("H7").Value = 14159
$ws.Range("I7").Value = 1590
$ws.Range("J7").Value = 119
$ws.Range("K7").Value = 132
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = "Bag"
